$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the marking scheme: per-right-answer mark goes from 5 to 4,
# and per-wrong-answer penalty goes from -1 to -2.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Recalculate totals based on corrected marking scheme.
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "58 / 112"
